$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same formatting as the rest of the header row
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every player row (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 80   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 82   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
